# Generate Report for Handoff
# Updates the localization-status workbook to mark the
# c7a5e5b8-d6f2-4c2b-81a8-4ab312ec7970 file as "Ready for handoff" across
# the Overview, zh-cn and de-de sheets, refreshes its handoff timestamps
# and priority, and widens the status columns to fit the new text.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet: row 3 is the c7a5e5b8-... file ----
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E3").Value = "Ready for handoff"
$ovw.Range("F3").Value = "Ready for handoff"
$ovw.Range("G3").Value = "2016-09-04 10:17:10"
$ovw.Columns.Item(5).ColumnWidth = 17.2159881591797
$ovw.Columns.Item(6).ColumnWidth = 17.2159881591797

# ---- zh-cn sheet: row 3 is the c7a5e5b8-... file ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("E3").Value = "mt"
$zh.Range("H3").Value = "2016-09-04 10:17:02"
$zh.Columns.Item(3).ColumnWidth = 17.2159881591797

# ---- de-de sheet: row 3 is the c7a5e5b8-... file ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Ready for handoff"
$de.Range("E3").Value = "mt"
$de.Range("H3").Value = "2016-09-04 10:17:10"
$de.Columns.Item(3).ColumnWidth = 17.2159881591797
